$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.679.11"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "3.145.26"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.40%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.577"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.19%  "
$ws.Range("D9").Value = "3.160.79"
$ws.Range("E9").Value = "  -0.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.118"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.14%  "
$ws.Range("E11").Value = "  -2.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.386"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("D13").Value = "3.692.15"
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.128"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").Value = "64.617.94"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.88%  "
$ws.Range("D17").Value = "3.145.67"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("E18").Value = "  -1.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "410.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.485"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.30%  "
$ws.Range("E26").Value = "  -5.08%  "
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "163.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.14"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.37%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("D39").Value = "2.644.02"
$ws.Range("E39").Value = "  -1.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.91%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.693"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0614"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.49%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "291.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("E48").Value = "  -2.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.995"
$ws.Range("D49").Style = "Normal"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.68%  "
